$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Delete the two removed rows (descending order to keep indices stable) ---
# Row 28 = "SC 92"
$ws.Rows.Item(28).Delete()
# Row 26 = "RM 232"
$ws.Rows.Item(26).Delete()

# --- Fill in / clear individual cell values per the diff ---
# E3 (RM 8): empty -> -5.7
$ws.Range("E3").Value = -5.7

# F4 (RM 9): 17.97 -> empty
$ws.Range("F4").ClearContents()

# E5 (RM 14): -5 -> empty
$ws.Range("E5").ClearContents()

# F9 (RM 42): empty -> 17.26
$ws.Range("F9").Value = 17.26

# F10 (RM 52 a): empty -> 16.43
$ws.Range("F10").Value = 16.43

# F17 (RM 116): 17.78 -> empty
$ws.Range("F17").ClearContents()

# F18 (RM 120): 18.35 -> empty
$ws.Range("F18").ClearContents()

# E21 (RM 135): empty -> -8.699999999999999
$ws.Range("E21").Value = -8.699999999999999

# E23 (RM 140): -7 -> empty
$ws.Range("E23").ClearContents()

# E32 (SC 193, after row shift): empty -> -6.4
$ws.Range("E32").Value = -6.4
